$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) for the season-record columns ---------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold / bordered / centered) by
# copying the style of the adjacent "Unnamed: 28" header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Season record values for every roster row (2-42) ---------------------
$wins = 93
$losses = 69
$ties = 0

for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # AD
    $ws.Cells.Item($row, 31).Value = $losses  # AE
    $ws.Cells.Item($row, 32).Value = $ties    # AF
}
